# Fix invalid (unquoted) sheet-name references in VLOOKUP formulas
# across the workbook, matching the OOXML diff (quote sheet names in A1 refs).
$wb = $excel.ActiveWorkbook

# Sheet: Таблицы
$ws = $wb.Worksheets.Item('Таблицы')
$ws.Range('B2').Formula = '=VLOOKUP(A2,''Группы''!A:B,2,0)'
$ws.Range('B3').Formula = '=VLOOKUP(A3,''Группы''!A:B,2,0)'
$ws.Range('B4').Formula = '=VLOOKUP(A4,''Группы''!A:B,2,0)'
$ws.Range('B5').Formula = '=VLOOKUP(A5,''Группы''!A:B,2,0)'
$ws.Range('B6').Formula = '=VLOOKUP(A6,''Группы''!A:B,2,0)'
$ws.Range('B7').Formula = '=VLOOKUP(A7,''Группы''!A:B,2,0)'
$ws.Range('B8').Formula = '=VLOOKUP(A8,''Группы''!A:B,2,0)'
$ws.Range('B9').Formula = '=VLOOKUP(A9,''Группы''!A:B,2,0)'
$ws.Range('B10').Formula = '=VLOOKUP(A10,''Группы''!A:B,2,0)'
$ws.Range('B11').Formula = '=VLOOKUP(A11,''Группы''!A:B,2,0)'
$ws.Range('B12').Formula = '=VLOOKUP(A12,''Группы''!A:B,2,0)'

# Sheet: Поля таблиц
$ws = $wb.Worksheets.Item('Поля таблиц')
$ws.Range('C2').Formula = '=VLOOKUP(B2,''Таблицы''!C:D,2,0)'
$ws.Range('C3').Formula = '=VLOOKUP(B3,''Таблицы''!C:D,2,0)'
$ws.Range('C4').Formula = '=VLOOKUP(B4,''Таблицы''!C:D,2,0)'
$ws.Range('C5').Formula = '=VLOOKUP(B5,''Таблицы''!C:D,2,0)'
$ws.Range('C6').Formula = '=VLOOKUP(B6,''Таблицы''!C:D,2,0)'
$ws.Range('C7').Formula = '=VLOOKUP(B7,''Таблицы''!C:D,2,0)'
$ws.Range('C8').Formula = '=VLOOKUP(B8,''Таблицы''!C:D,2,0)'
$ws.Range('C9').Formula = '=VLOOKUP(B9,''Таблицы''!C:D,2,0)'
$ws.Range('C10').Formula = '=VLOOKUP(B10,''Таблицы''!C:D,2,0)'
$ws.Range('C11').Formula = '=VLOOKUP(B11,''Таблицы''!C:D,2,0)'
$ws.Range('C12').Formula = '=VLOOKUP(B12,''Таблицы''!C:D,2,0)'
$ws.Range('C13').Formula = '=VLOOKUP(B13,''Таблицы''!C:D,2,0)'
$ws.Range('C14').Formula = '=VLOOKUP(B14,''Таблицы''!C:D,2,0)'
$ws.Range('C15').Formula = '=VLOOKUP(B15,''Таблицы''!C:D,2,0)'
$ws.Range('C16').Formula = '=VLOOKUP(B16,''Таблицы''!C:D,2,0)'
$ws.Range('C17').Formula = '=VLOOKUP(B17,''Таблицы''!C:D,2,0)'
$ws.Range('C18').Formula = '=VLOOKUP(B18,''Таблицы''!C:D,2,0)'
$ws.Range('C19').Formula = '=VLOOKUP(B19,''Таблицы''!C:D,2,0)'
$ws.Range('C20').Formula = '=VLOOKUP(B20,''Таблицы''!C:D,2,0)'
$ws.Range('C21').Formula = '=VLOOKUP(B21,''Таблицы''!C:D,2,0)'
$ws.Range('C22').Formula = '=VLOOKUP(B22,''Таблицы''!C:D,2,0)'
$ws.Range('C23').Formula = '=VLOOKUP(B23,''Таблицы''!C:D,2,0)'
$ws.Range('C24').Formula = '=VLOOKUP(B24,''Таблицы''!C:D,2,0)'
$ws.Range('C25').Formula = '=VLOOKUP(B25,''Таблицы''!C:D,2,0)'
$ws.Range('C26').Formula = '=VLOOKUP(B26,''Таблицы''!C:D,2,0)'
$ws.Range('C27').Formula = '=VLOOKUP(B27,''Таблицы''!C:D,2,0)'
$ws.Range('C28').Formula = '=VLOOKUP(B28,''Таблицы''!C:D,2,0)'
$ws.Range('C29').Formula = '=VLOOKUP(B29,''Таблицы''!C:D,2,0)'
$ws.Range('C30').Formula = '=VLOOKUP(B30,''Таблицы''!C:D,2,0)'
$ws.Range('C31').Formula = '=VLOOKUP(B31,''Таблицы''!C:D,2,0)'
$ws.Range('C32').Formula = '=VLOOKUP(B32,''Таблицы''!C:D,2,0)'
$ws.Range('C33').Formula = '=VLOOKUP(B33,''Таблицы''!C:D,2,0)'
$ws.Range('C34').Formula = '=VLOOKUP(B34,''Таблицы''!C:D,2,0)'
$ws.Range('C35').Formula = '=VLOOKUP(B35,''Таблицы''!C:D,2,0)'
$ws.Range('C36').Formula = '=VLOOKUP(B36,''Таблицы''!C:D,2,0)'
$ws.Range('C37').Formula = '=VLOOKUP(B37,''Таблицы''!C:D,2,0)'
$ws.Range('C38').Formula = '=VLOOKUP(B38,''Таблицы''!C:D,2,0)'
$ws.Range('C39').Formula = '=VLOOKUP(B39,''Таблицы''!C:D,2,0)'
$ws.Range('C40').Formula = '=VLOOKUP(B40,''Таблицы''!C:D,2,0)'
$ws.Range('C41').Formula = '=VLOOKUP(B41,''Таблицы''!C:D,2,0)'
$ws.Range('C42').Formula = '=VLOOKUP(B42,''Таблицы''!C:D,2,0)'
$ws.Range('C43').Formula = '=VLOOKUP(B43,''Таблицы''!C:D,2,0)'
$ws.Range('C44').Formula = '=VLOOKUP(B44,''Таблицы''!C:D,2,0)'
$ws.Range('C45').Formula = '=VLOOKUP(B45,''Таблицы''!C:D,2,0)'
$ws.Range('C46').Formula = '=VLOOKUP(B46,''Таблицы''!C:D,2,0)'
$ws.Range('C47').Formula = '=VLOOKUP(B47,''Таблицы''!C:D,2,0)'
$ws.Range('C48').Formula = '=VLOOKUP(B48,''Таблицы''!C:D,2,0)'
$ws.Range('C49').Formula = '=VLOOKUP(B49,''Таблицы''!C:D,2,0)'
$ws.Range('C50').Formula = '=VLOOKUP(B50,''Таблицы''!C:D,2,0)'
$ws.Range('C51').Formula = '=VLOOKUP(B51,''Таблицы''!C:D,2,0)'
$ws.Range('C52').Formula = '=VLOOKUP(B52,''Таблицы''!C:D,2,0)'
$ws.Range('C53').Formula = '=VLOOKUP(B53,''Таблицы''!C:D,2,0)'
$ws.Range('C54').Formula = '=VLOOKUP(B54,''Таблицы''!C:D,2,0)'
$ws.Range('C55').Formula = '=VLOOKUP(B55,''Таблицы''!C:D,2,0)'

# Sheet: Индексы
$ws = $wb.Worksheets.Item('Индексы')
$ws.Range('C2').Formula = '=VLOOKUP(B2,''Таблицы''!C:D,2,0)'

# Sheet: Поля индексов
$ws = $wb.Worksheets.Item('Поля индексов')
$ws.Range('B2').Formula = '=VLOOKUP(A2,''Индексы''!A:E,2,0)'
$ws.Range('C2').Formula = '=VLOOKUP(A2,''Индексы''!A:E,3,0)'
$ws.Range('D2').Formula = '=VLOOKUP(A2,''Индексы''!A:E,4,0)'
$ws.Range('E2').Formula = '=VLOOKUP(A2,''Индексы''!A:E,5,0)'
$ws.Range('G2').Formula = '=VLOOKUP(B2&"."&F2,''Поля таблиц''!A:G,7,0)'
$ws.Range('B3').Formula = '=VLOOKUP(A3,''Индексы''!A:E,2,0)'
$ws.Range('C3').Formula = '=VLOOKUP(A3,''Индексы''!A:E,3,0)'
$ws.Range('D3').Formula = '=VLOOKUP(A3,''Индексы''!A:E,4,0)'
$ws.Range('E3').Formula = '=VLOOKUP(A3,''Индексы''!A:E,5,0)'
$ws.Range('G3').Formula = '=VLOOKUP(B3&"."&F3,''Поля таблиц''!A:G,7,0)'

# Sheet: Отношения
$ws = $wb.Worksheets.Item('Отношения')
$ws.Range('C2').Formula = '=VLOOKUP(B2,''Таблицы''!C:D,2,0)'
$ws.Range('H2').Formula = '=VLOOKUP(G2,''Таблицы''!C:D,2,0)'
$ws.Range('C3').Formula = '=VLOOKUP(B3,''Таблицы''!C:D,2,0)'
$ws.Range('H3').Formula = '=VLOOKUP(G3,''Таблицы''!C:D,2,0)'
$ws.Range('C4').Formula = '=VLOOKUP(B4,''Таблицы''!C:D,2,0)'
$ws.Range('H4').Formula = '=VLOOKUP(G4,''Таблицы''!C:D,2,0)'
$ws.Range('C5').Formula = '=VLOOKUP(B5,''Таблицы''!C:D,2,0)'
$ws.Range('H5').Formula = '=VLOOKUP(G5,''Таблицы''!C:D,2,0)'
$ws.Range('C6').Formula = '=VLOOKUP(B6,''Таблицы''!C:D,2,0)'
$ws.Range('H6').Formula = '=VLOOKUP(G6,''Таблицы''!C:D,2,0)'
$ws.Range('C7').Formula = '=VLOOKUP(B7,''Таблицы''!C:D,2,0)'
$ws.Range('H7').Formula = '=VLOOKUP(G7,''Таблицы''!C:D,2,0)'
$ws.Range('C8').Formula = '=VLOOKUP(B8,''Таблицы''!C:D,2,0)'
$ws.Range('H8').Formula = '=VLOOKUP(G8,''Таблицы''!C:D,2,0)'
$ws.Range('C9').Formula = '=VLOOKUP(B9,''Таблицы''!C:D,2,0)'
$ws.Range('H9').Formula = '=VLOOKUP(G9,''Таблицы''!C:D,2,0)'

# Sheet: Поля отношений
$ws = $wb.Worksheets.Item('Поля отношений')
$ws.Range('B2').Formula = '=VLOOKUP(A2,''Отношения''!A:E,2,0)'
$ws.Range('C2').Formula = '=VLOOKUP(A2,''Отношения''!A:E,3,0)'
$ws.Range('E2').Formula = '=VLOOKUP(D2,''Таблицы''!C:D,2,0)'
$ws.Range('F2').Formula = '=VLOOKUP(A2,''Отношения''!A:E,4,0)'
$ws.Range('G2').Formula = '=VLOOKUP(A2,''Отношения''!A:E,5,0)'
$ws.Range('I2').Formula = '=VLOOKUP(B2&"."&H2,''Поля таблиц''!A:G,7,0)'
$ws.Range('K2').Formula = '=VLOOKUP(D2&"."&J2,''Поля таблиц''!A:G,7,0)'
$ws.Range('B3').Formula = '=VLOOKUP(A3,''Отношения''!A:E,2,0)'
$ws.Range('C3').Formula = '=VLOOKUP(A3,''Отношения''!A:E,3,0)'
$ws.Range('E3').Formula = '=VLOOKUP(D3,''Таблицы''!C:D,2,0)'
$ws.Range('F3').Formula = '=VLOOKUP(A3,''Отношения''!A:E,4,0)'
$ws.Range('G3').Formula = '=VLOOKUP(A3,''Отношения''!A:E,5,0)'
$ws.Range('I3').Formula = '=VLOOKUP(B3&"."&H3,''Поля таблиц''!A:G,7,0)'
$ws.Range('K3').Formula = '=VLOOKUP(D3&"."&J3,''Поля таблиц''!A:G,7,0)'
$ws.Range('B4').Formula = '=VLOOKUP(A4,''Отношения''!A:E,2,0)'
$ws.Range('C4').Formula = '=VLOOKUP(A4,''Отношения''!A:E,3,0)'
$ws.Range('E4').Formula = '=VLOOKUP(D4,''Таблицы''!C:D,2,0)'
$ws.Range('F4').Formula = '=VLOOKUP(A4,''Отношения''!A:E,4,0)'
$ws.Range('G4').Formula = '=VLOOKUP(A4,''Отношения''!A:E,5,0)'
$ws.Range('I4').Formula = '=VLOOKUP(B4&"."&H4,''Поля таблиц''!A:G,7,0)'
$ws.Range('K4').Formula = '=VLOOKUP(D4&"."&J4,''Поля таблиц''!A:G,7,0)'
$ws.Range('B5').Formula = '=VLOOKUP(A5,''Отношения''!A:E,2,0)'
$ws.Range('C5').Formula = '=VLOOKUP(A5,''Отношения''!A:E,3,0)'
$ws.Range('E5').Formula = '=VLOOKUP(D5,''Таблицы''!C:D,2,0)'
$ws.Range('F5').Formula = '=VLOOKUP(A5,''Отношения''!A:E,4,0)'
$ws.Range('G5').Formula = '=VLOOKUP(A5,''Отношения''!A:E,5,0)'
$ws.Range('I5').Formula = '=VLOOKUP(B5&"."&H5,''Поля таблиц''!A:G,7,0)'
$ws.Range('K5').Formula = '=VLOOKUP(D5&"."&J5,''Поля таблиц''!A:G,7,0)'
$ws.Range('B6').Formula = '=VLOOKUP(A6,''Отношения''!A:E,2,0)'
$ws.Range('C6').Formula = '=VLOOKUP(A6,''Отношения''!A:E,3,0)'
$ws.Range('E6').Formula = '=VLOOKUP(D6,''Таблицы''!C:D,2,0)'
$ws.Range('F6').Formula = '=VLOOKUP(A6,''Отношения''!A:E,4,0)'
$ws.Range('G6').Formula = '=VLOOKUP(A6,''Отношения''!A:E,5,0)'
$ws.Range('I6').Formula = '=VLOOKUP(B6&"."&H6,''Поля таблиц''!A:G,7,0)'
$ws.Range('K6').Formula = '=VLOOKUP(D6&"."&J6,''Поля таблиц''!A:G,7,0)'
$ws.Range('B7').Formula = '=VLOOKUP(A7,''Отношения''!A:E,2,0)'
$ws.Range('C7').Formula = '=VLOOKUP(A7,''Отношения''!A:E,3,0)'
$ws.Range('E7').Formula = '=VLOOKUP(D7,''Таблицы''!C:D,2,0)'
$ws.Range('F7').Formula = '=VLOOKUP(A7,''Отношения''!A:E,4,0)'
$ws.Range('G7').Formula = '=VLOOKUP(A7,''Отношения''!A:E,5,0)'
$ws.Range('I7').Formula = '=VLOOKUP(B7&"."&H7,''Поля таблиц''!A:G,7,0)'
$ws.Range('K7').Formula = '=VLOOKUP(D7&"."&J7,''Поля таблиц''!A:G,7,0)'
$ws.Range('B8').Formula = '=VLOOKUP(A8,''Отношения''!A:E,2,0)'
$ws.Range('C8').Formula = '=VLOOKUP(A8,''Отношения''!A:E,3,0)'
$ws.Range('E8').Formula = '=VLOOKUP(D8,''Таблицы''!C:D,2,0)'
$ws.Range('F8').Formula = '=VLOOKUP(A8,''Отношения''!A:E,4,0)'
$ws.Range('G8').Formula = '=VLOOKUP(A8,''Отношения''!A:E,5,0)'
$ws.Range('I8').Formula = '=VLOOKUP(B8&"."&H8,''Поля таблиц''!A:G,7,0)'
$ws.Range('K8').Formula = '=VLOOKUP(D8&"."&J8,''Поля таблиц''!A:G,7,0)'
$ws.Range('B9').Formula = '=VLOOKUP(A9,''Отношения''!A:E,2,0)'
$ws.Range('C9').Formula = '=VLOOKUP(A9,''Отношения''!A:E,3,0)'
$ws.Range('E9').Formula = '=VLOOKUP(D9,''Таблицы''!C:D,2,0)'
$ws.Range('F9').Formula = '=VLOOKUP(A9,''Отношения''!A:E,4,0)'
$ws.Range('G9').Formula = '=VLOOKUP(A9,''Отношения''!A:E,5,0)'
$ws.Range('I9').Formula = '=VLOOKUP(B9&"."&H9,''Поля таблиц''!A:G,7,0)'
$ws.Range('K9').Formula = '=VLOOKUP(D9&"."&J9,''Поля таблиц''!A:G,7,0)'
$ws.Range('B10').Formula = '=VLOOKUP(A10,''Отношения''!A:E,2,0)'
$ws.Range('C10').Formula = '=VLOOKUP(A10,''Отношения''!A:E,3,0)'
$ws.Range('E10').Formula = '=VLOOKUP(D10,''Таблицы''!C:D,2,0)'
$ws.Range('F10').Formula = '=VLOOKUP(A10,''Отношения''!A:E,4,0)'
$ws.Range('G10').Formula = '=VLOOKUP(A10,''Отношения''!A:E,5,0)'
$ws.Range('I10').Formula = '=VLOOKUP(B10&"."&H10,''Поля таблиц''!A:G,7,0)'
$ws.Range('K10').Formula = '=VLOOKUP(D10&"."&J10,''Поля таблиц''!A:G,7,0)'

# Sheet: Элементы перечислений
$ws = $wb.Worksheets.Item('Элементы перечислений')
$ws.Range('B2').Formula = '=VLOOKUP(A2,''Перечисления''!A:B,2,0)'
$ws.Range('B3').Formula = '=VLOOKUP(A3,''Перечисления''!A:B,2,0)'
$ws.Range('B4').Formula = '=VLOOKUP(A4,''Перечисления''!A:B,2,0)'
